$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing header cell L2: "Stations" -> "0.15-0.25 Stations" ---
$ws.Range("L2").Value = "0.15-0.25 Stations"

# --- Copy the existing header style (B2:L2, centered + border) onto the new header cells M2:Q2 ---
$ws.Range("L2").Copy() | Out-Null
$ws.Range("M2:Q2").PasteSpecial(-4122) | Out-Null

# New sub-headers for the "0.15-0.25" group (M2:N2) and "17-18" group (O2:Q2)
$ws.Range("M2").Value = "0.15-0.25 Latitude"
$ws.Range("N2").Value = "0.15-0.25 Times"
$ws.Range("O2").Value = "17-18° Stations"
$ws.Range("P2").Value = "17-18° Latitude"
$ws.Range("Q2").Value = "17-18° Times"

# --- Data rows (2008, 2009, 2011) for the new "Stations" / "Latitude" / "Times" columns ---

# Row 3 (2008)
$ws.Range("L3").Value = "35-41"
$ws.Range("M3").Value = "34.25-35.75"
$ws.Range("N3").Value = "3/29-4/3"
$ws.Range("O3").Value = "27-30"
$ws.Range("P3").Value = "32.25-33"
$ws.Range("Q3").Value = "3/28-3/29"

# Row 4 (2009)
$ws.Range("L4").Value = "33-35"
$ws.Range("M4").Value = "35-35.75"
$ws.Range("N4").Value = "3/22-3/24"
$ws.Range("O4").Value = "18-22"
$ws.Range("P4").Value = "31.25-32.25"
$ws.Range("Q4").Value = 41719
$ws.Range("N4").NumberFormat = "d-mmm"
$ws.Range("Q4").NumberFormat = "d-mmm"

# Row 5 (2011)
$ws.Range("L5").Value = "26-30"
$ws.Range("M5").Value = "32.5-33.5"
$ws.Range("N5").Value = "3/19-3/20"
$ws.Range("O5").Value = "21-25"
$ws.Range("P5").Value = "31.25-32.25"
$ws.Range("Q5").Value = 41718
$ws.Range("Q5").NumberFormat = "d-mmm"

# --- Column widths for the new columns (approximate character widths used by the source workbook) ---
$ws.Range("L1").ColumnWidth = 15.333333333333334
$ws.Range("M1").ColumnWidth = 15.333333333333334
$ws.Range("N1").ColumnWidth = 15.333333333333334
$ws.Range("O1").ColumnWidth = 13.0
$ws.Range("P1").ColumnWidth = 13.0
$ws.Range("Q1").ColumnWidth = 13.5

# --- Restore view/selection state ---
$ws.Range("N5").Select() | Out-Null

Write-Output "Edit complete"
